# Apply "new progress as of date 04 nov 2025" update
# For rows 3-20 on the "Training Dashboard" sheet:
#   - Column H (PERIOD TO EXPIRE) decreases by 1
#   - Column I (LAST UPDATE) changes from 03-Nov-2025 to 04-Nov-2025

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Training Dashboard")

# Column I holds a date-like string (e.g. "03-Nov-2025") stored as TEXT.
# Mark the range as Text first so the new date string isn't auto-converted
# into a real date serial when it is written back.
$ws.Range("I3:I20").NumberFormat = "@"

for ($row = 3; $row -le 20; $row++) {
    $hCell = $ws.Cells.Item($row, 8)   # Column H
    $iCell = $ws.Cells.Item($row, 9)   # Column I

    $hVal = $hCell.Value2
    $hCell.Value2 = $hVal - 1

    $iCell.Value2 = "04-Nov-2025"
}
